# thesis gantt chart.xlsx - apply the author's edits:
#  - advance the Gantt chart "current week" scroll value from 9 to 15
#    (this cascades through all the Week-N / date header formulas automatically)
#  - update several task duration (G) / percent-complete (H) values
#  - clear a couple of stray cell values
#  - update the frozen-pane scroll position / active selection

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GanttChart")

# --- Main "current week" control cell -------------------------------------
# Drives K6/K5/... and all the "Week N" header formulas via I4.
$ws.Range("I4").Value = 15

# --- Task table updates (columns G = duration(days), H = % complete) ------

# Row 23 - Task 2 "Proposal" duration 80 -> 77 days
$ws.Range("G23").Value = 77

# Row 27 - Task 2.4, % complete 0.7 -> 1 (100%)
$ws.Range("H27").Value = 1

# Row 28 - Task 2.5
$ws.Range("C28").ClearContents()
$ws.Range("G28").Value = 18
$ws.Range("H28").Value = 1

# Row 31 - Task 2.8
$ws.Range("G31").Value = 20
$ws.Range("H31").Value = 0.05

# Row 34
$ws.Range("G34").Value = 8

# Row 35 - fully clear the WBS id cell (contents + formatting)
$ws.Range("A35").Clear()

# Row 36 - the "end" milestone date is pasted as a fixed value instead of
# the old "=F23" formula
$ws.Range("E36").Value = 45429

# --- View state: scroll the frozen pane down and move the selection -------
$ws.Range("AW25").Select()
